$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the 11-year tax-year column headers forward by one year
# (2021-2031 -> 2022-2032) to roll the table onto the May 2022 baseline.
$ws.Range("B9").Value = 2022
$ws.Range("C9").Value = 2023
$ws.Range("D9").Value = 2024
$ws.Range("E9").Value = 2025
$ws.Range("F9").Value = 2026
$ws.Range("G9").Value = 2027
$ws.Range("H9").Value = 2028
$ws.Range("I9").Value = 2029
$ws.Range("J9").Value = 2030
$ws.Range("K9").Value = 2031
$ws.Range("L9").Value = 2032

# Tax Wedge for Business Assets
$ws.Range("B11").Value = 0.98
$ws.Range("C11").Value = 1.02
$ws.Range("D11").Value = 1.07
$ws.Range("E11").Value = 1.11
$ws.Range("F11").Value = 1.23
$ws.Range("G11").Value = 1.28
$ws.Range("H11").Value = 1.28
$ws.Range("I11").Value = 1.28
$ws.Range("J11").Value = 1.29
$ws.Range("K11").Value = 1.29
$ws.Range("L11").Value = 1.3

# By source of financing: Equity-financed
$ws.Range("B14").Value = 1.2
$ws.Range("C14").Value = 1.25
$ws.Range("D14").Value = 1.3
$ws.Range("E14").Value = 1.35
$ws.Range("F14").Value = 1.52
$ws.Range("G14").Value = 1.57
$ws.Range("H14").Value = 1.57
$ws.Range("I14").Value = 1.57
$ws.Range("J14").Value = 1.58
$ws.Range("K14").Value = 1.57
$ws.Range("L14").Value = 1.58

# By source of financing: Debt-financed
$ws.Range("B15").Value = 0.53
$ws.Range("C15").Value = 0.55
$ws.Range("D15").Value = 0.59
$ws.Range("E15").Value = 0.62
$ws.Range("F15").Value = 0.62
$ws.Range("G15").Value = 0.67
$ws.Range("H15").Value = 0.68
$ws.Range("I15").Value = 0.68
$ws.Range("J15").Value = 0.6899999999999999
$ws.Range("K15").Value = 0.7
$ws.Range("L15").Value = 0.71

# Difference between sources of financing
$ws.Range("B17").Value = 0.67
$ws.Range("C17").Value = 0.6899999999999999
$ws.Range("D17").Value = 0.71
$ws.Range("E17").Value = 0.73
$ws.Range("F17").Value = 0.89
$ws.Range("G17").Value = 0.9
$ws.Range("H17").Value = 0.89
$ws.Range("I17").Value = 0.89
$ws.Range("J17").Value = 0.88
$ws.Range("K17").Value = 0.88
$ws.Range("L17").Value = 0.87

# By legal form of organization: C corporations
$ws.Range("B20").Value = 1.08
$ws.Range("C20").Value = 1.12
$ws.Range("D20").Value = 1.16
$ws.Range("E20").Value = 1.2
$ws.Range("F20").Value = 1.25
$ws.Range("G20").Value = 1.3
$ws.Range("H20").Value = 1.3
$ws.Range("I20").Value = 1.3
$ws.Range("J20").Value = 1.31
$ws.Range("K20").Value = 1.31
$ws.Range("L20").Value = 1.32

# By legal form of organization: Pass-through entities
$ws.Range("B21").Value = 0.79
$ws.Range("C21").Value = 0.84
$ws.Range("D21").Value = 0.89
$ws.Range("E21").Value = 0.9399999999999999
$ws.Range("F21").Value = 1.17
$ws.Range("G21").Value = 1.24
$ws.Range("H21").Value = 1.24
$ws.Range("I21").Value = 1.24
$ws.Range("J21").Value = 1.25
$ws.Range("K21").Value = 1.25
$ws.Range("L21").Value = 1.27

# Difference between legal forms of organization
$ws.Range("B23").Value = 0.3
$ws.Range("C23").Value = 0.28
$ws.Range("D23").Value = 0.27
$ws.Range("E23").Value = 0.26
$ws.Range("F23").Value = 0.08
$ws.Range("G23").Value = 0.06
$ws.Range("H23").Value = 0.06
$ws.Range("I23").Value = 0.06
$ws.Range("J23").Value = 0.06
$ws.Range("K23").Value = 0.05
$ws.Range("L23").Value = 0.05

# Weighted mean absolute difference between all asset pairs
$ws.Range("B25").Value = 0.87
$ws.Range("C25").Value = 0.8100000000000001
$ws.Range("D25").Value = 0.76
$ws.Range("E25").Value = 0.71
$ws.Range("F25").Value = 0.7
$ws.Range("G25").Value = 0.66
$ws.Range("H25").Value = 0.66
$ws.Range("I25").Value = 0.66
$ws.Range("J25").Value = 0.66
$ws.Range("K25").Value = 0.66
$ws.Range("L25").Value = 0.66

# Weighted mean absolute difference between all industry pairs
$ws.Range("B26").Value = 0.06
$ws.Range("D26").Value = 0.04

# Tax Wedge for Owner-Occupied Housing
$ws.Range("B28").Value = 0.05
$ws.Range("C28").Value = 0.05
$ws.Range("D28").Value = 0.02
$ws.Range("E28").Value = -0.02
$ws.Range("F28").Value = -0.5
$ws.Range("G28").Value = -0.5
$ws.Range("H28").Value = -0.5
$ws.Range("I28").Value = -0.5
$ws.Range("L28").Value = -0.57

# Difference between owner-occupied housing and business assets
$ws.Range("B29").Value = -0.93
$ws.Range("C29").Value = -0.97
$ws.Range("D29").Value = -1.05
$ws.Range("E29").Value = -1.13
$ws.Range("F29").Value = -1.72
$ws.Range("G29").Value = -1.78
$ws.Range("H29").Value = -1.78
$ws.Range("I29").Value = -1.78
$ws.Range("J29").Value = -1.79
$ws.Range("K29").Value = -1.79
$ws.Range("L29").Value = -1.87

# By source of financing (owner-occupied housing): Equity-financed
$ws.Range("F32").Value = -0.25
$ws.Range("H32").Value = -0.25
$ws.Range("I32").Value = -0.25
$ws.Range("J32").Value = -0.25
$ws.Range("K32").Value = -0.25
$ws.Range("L32").Value = -0.32

# By source of financing (owner-occupied housing): Debt-financed
$ws.Range("B33").Value = 0.19
$ws.Range("C33").Value = 0.17
$ws.Range("D33").Value = 0.09
$ws.Range("E33").Value = -0.02
$ws.Range("F33").Value = -0.93
$ws.Range("G33").Value = -0.93
$ws.Range("H33").Value = -0.93
$ws.Range("I33").Value = -0.9399999999999999
$ws.Range("J33").Value = -0.9399999999999999
$ws.Range("K33").Value = -0.95
$ws.Range("L33").Value = -1.02

# Difference between sources of financing (owner-occupied housing)
$ws.Range("B35").Value = -0.2
$ws.Range("C35").Value = -0.19
$ws.Range("D35").Value = -0.11
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0.68
$ws.Range("G35").Value = 0.68
$ws.Range("H35").Value = 0.68
$ws.Range("I35").Value = 0.6899999999999999
$ws.Range("J35").Value = 0.6899999999999999
$ws.Range("K35").Value = 0.6899999999999999
$ws.Range("L35").Value = 0.7

# Memorandum: Tax Wedge for Owner-Occupied Housing Structures
$ws.Range("B38").Value = 0.05
$ws.Range("C38").Value = 0.05
$ws.Range("D38").Value = 0.02
$ws.Range("E38").Value = -0.02
$ws.Range("F38").Value = -0.5
$ws.Range("G38").Value = -0.5
$ws.Range("H38").Value = -0.5
$ws.Range("I38").Value = -0.5
$ws.Range("L38").Value = -0.57

# Memorandum: Tax Wedge for Renter-Occupied Housing Structures
$ws.Range("B39").Value = 1.4
$ws.Range("C39").Value = 1.39
$ws.Range("D39").Value = 1.38
$ws.Range("E39").Value = 1.37
$ws.Range("F39").Value = 1.45
$ws.Range("G39").Value = 1.45
$ws.Range("H39").Value = 1.45
$ws.Range("I39").Value = 1.45
$ws.Range("J39").Value = 1.46
$ws.Range("K39").Value = 1.46
$ws.Range("L39").Value = 1.47

# Memorandum: Difference between owner- and renter-occupied housing structures
$ws.Range("B41").Value = -1.35
$ws.Range("C41").Value = -1.34
$ws.Range("D41").Value = -1.36
$ws.Range("E41").Value = -1.39
$ws.Range("F41").Value = -1.95
$ws.Range("G41").Value = -1.94
$ws.Range("H41").Value = -1.95
$ws.Range("I41").Value = -1.95
$ws.Range("J41").Value = -1.96
$ws.Range("K41").Value = -1.96
$ws.Range("L41").Value = -2.04
